$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the poke arm length (L1) value
$ws.Range("B2").Value = 300

# Clear L2 value (C2) entirely
$ws.Range("C2").ClearContents()

# Update selection to reflect the new active cell
$ws.Range("X10").Select()

$wb.Save()
